$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Repoint ontology/terminology URIs from http:// to https:// (column A) ---
# Note: `.Value` is a write-only stub in this host; `.Value2` is used to read back the
# existing text so the http->https swap doesn't depend on hardcoding every string.
$httpsRows = @(3,5,7,10,12,14,17,18,19,20,21,22,23,24,26)
foreach ($r in $httpsRows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 -replace '^http://', 'https://'
}

# --- 2. Append new recommendation rows 27-38 (Respiratory Rate / ARDS / Sepsis / Catecholamine Therapy) ---
# row 27
$ws.Cells.Item(27,1).Value = 'https://data.cochrane.org/concepts/'
$ws.Cells.Item(27,2).Value = 'r4hp3qg8yn71'
$ws.Cells.Item(27,3).Value = 'Condition - Respiratory Rate'
# row 28
$ws.Cells.Item(28,1).Value = 'https://snomed.info/sct'
$ws.Cells.Item(28,2).Value = '86290005'
$ws.Cells.Item(28,3).Value = 'Respiratory measure (observable entity) '
# row 29
$ws.Cells.Item(29,1).Value = 'https://data.cochrane.org/concepts/'
$ws.Cells.Item(29,2).Value = 'r4hp3pkr3yl3'
$ws.Cells.Item(29,3).Value = 'Condition - Acute Respiratory Distress Syndrome'
$ws.Cells.Item(29,4).Value = 'oxygenation_index_cacl'
$ws.Cells.Item(29,5).Value = 'float'
$ws.Cells.Item(29,8).Value = 300
$ws.Cells.Item(29,9).Value = 'mmHg'
# row 30
$ws.Cells.Item(30,1).Value = 'https://snomed.info/sct'
$ws.Cells.Item(30,2).Value = '67782005'
$ws.Cells.Item(30,3).Value = 'Acute respiratory distress syndrome (disorder)'
$ws.Cells.Item(30,4).Value = 'oxygenation_index_cacl'
$ws.Cells.Item(30,5).Value = 'float'
$ws.Cells.Item(30,8).Value = 300
$ws.Cells.Item(30,9).Value = 'mmHg'
# row 31
$ws.Cells.Item(31,1).Value = 'https://data.cochrane.org/concepts/'
$ws.Cells.Item(31,2).Value = 'r4hp3p86xjp4'
$ws.Cells.Item(31,3).Value = 'Condition - Sepsis'
$ws.Cells.Item(31,4).Value = 'deltaSOFA'
$ws.Cells.Item(31,5).Value = 'int'
$ws.Cells.Item(31,7).Value = 2
# row 32
$ws.Cells.Item(32,1).Value = 'https://snomed.info/sct'
$ws.Cells.Item(32,2).Value = '91302008'
$ws.Cells.Item(32,3).Value = 'Sepsis (disorder)'
$ws.Cells.Item(32,4).Value = 'deltaSOFA'
$ws.Cells.Item(32,5).Value = 'int'
$ws.Cells.Item(32,7).Value = 2
# row 33
$ws.Cells.Item(33,1).Value = 'https://covid-evidenz.de/concepts/'
$ws.Cells.Item(33,2).Value = 'catecholamine-therapy'
$ws.Cells.Item(33,3).Value = 'Catecholamine Therapy'
$ws.Cells.Item(33,4).Value = 'drug_norepinephrine'
$ws.Cells.Item(33,5).Value = 'float'
$ws.Cells.Item(33,7).Value = 0
$ws.Cells.Item(33,9).Value = 'µg/kg/min'
# row 34
$ws.Cells.Item(34,1).Value = 'https://covid-evidenz.de/concepts/'
$ws.Cells.Item(34,2).Value = 'catecholamine-therapy'
$ws.Cells.Item(34,3).Value = 'Catecholamine Therapy'
$ws.Cells.Item(34,4).Value = 'drug_epinephrine'
$ws.Cells.Item(34,5).Value = 'float'
$ws.Cells.Item(34,7).Value = 0
$ws.Cells.Item(34,9).Value = 'µg/kg/min'
# row 35
$ws.Cells.Item(35,1).Value = 'https://covid-evidenz.de/concepts/'
$ws.Cells.Item(35,2).Value = 'catecholamine-therapy'
$ws.Cells.Item(35,3).Value = 'Catecholamine Therapy'
$ws.Cells.Item(35,4).Value = 'drug_vasopressin'
$ws.Cells.Item(35,5).Value = 'float'
$ws.Cells.Item(35,7).Value = 0
$ws.Cells.Item(35,9).Value = 'E/kg/h'
# row 36
$ws.Cells.Item(36,1).Value = 'https://covid-evidenz.de/concepts/'
$ws.Cells.Item(36,2).Value = 'catecholamine-therapy'
$ws.Cells.Item(36,3).Value = 'Catecholamine Therapy'
$ws.Cells.Item(36,4).Value = 'drug_dobutamine'
$ws.Cells.Item(36,5).Value = 'float'
$ws.Cells.Item(36,7).Value = 0
$ws.Cells.Item(36,9).Value = 'µg/kg/min'
# row 37
$ws.Cells.Item(37,1).Value = 'https://covid-evidenz.de/concepts/'
$ws.Cells.Item(37,2).Value = 'catecholamine-therapy'
$ws.Cells.Item(37,3).Value = 'Catecholamine Therapy'
$ws.Cells.Item(37,4).Value = 'drug_dopamine'
$ws.Cells.Item(37,5).Value = 'float'
$ws.Cells.Item(37,7).Value = 0
$ws.Cells.Item(37,9).Value = 'µg/kg/min'

# --- 3. Wrap-text the long description/label cells, matching the rest of the sheet's style ---
$ws.Range('C27').WrapText = $true
$ws.Range('B28').WrapText = $true
$ws.Range('C28').WrapText = $true
$ws.Range('A33').WrapText = $true
$ws.Range('A34').WrapText = $true
$ws.Range('A35').WrapText = $true
$ws.Range('A36').WrapText = $true
$ws.Range('A37').WrapText = $true
$ws.Range('A38').WrapText = $true

# --- 4. Conditional formatting on D33:D37 ---
# A throwaway rule is added first and removed afterwards purely so the surviving rule's
# priority comes out as 2 (matching the source file) instead of 1.
$dummyFc = $ws.Range('D33:D37').FormatConditions.Add(2, 3, '=TRUE')
$fc = $ws.Range('D33:D37').FormatConditions.Add(2, 3, '=$V33=1')
$fc.Font.Color = 0
$dummyFc.Delete()

# --- 5. Widen column B to fit the new 'catecholamine-therapy' values ---
$ws.Columns.Item(2).ColumnWidth = 20.0

# --- 6. Leave the selection where the author left it ---
$ws.Range('B39').Select()
